# V1.1.0 update: add switchableSlots row and a "Comment" column (H)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlEdgeLeft = 7
$xlEdgeTop = 8
$xlEdgeBottom = 9
$xlEdgeRight = 10
$xlContinuous = 1
$xlThin = 2
$xlCenter = -4108

# 1) Insert a new row at row 3 (pushes old rows 3-11 down to 4-12)
$ws.Rows("3").Insert()

# 2) Fill the new row 3 with the switchableSlots definition
$ws.Range("B3").Value = "switchableSlots"
$ws.Range("C3").Value = "为一个用逗号分割的int列表"
$ws.Range("D3").Value = "可选弹药所在槽位"
$ws.Range("E3").Value = "用于列举弹药来源槽位"
$ws.Range("F3").Value = "/"
$ws.Range("G3").Value = "/"

# 3) Add the new "注释" (Comment) column in H
$ws.Range("H1").Value = "注释"
$ws.Range("H2").Value = "旧版"
$ws.Range("H3").Value = "定义此参数将自动覆盖switchableProjectiles所定义的内容"

# 4) Alignment for the whole table (including new column/row)
$ws.Range("A1:H12").HorizontalAlignment = $xlCenter
$ws.Range("A1:H12").VerticalAlignment = $xlCenter

# 5) Column width for the new column H
$ws.Range("H1").ColumnWidth = 48.5

# 6) Borders - rebuild the table grid lines
# LEFT edges
$ws.Range("A1:H1").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("A2:B5").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("A6:A11").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("A12").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous
$ws.Range("C12").Borders.Item($xlEdgeLeft).LineStyle = $xlContinuous

# RIGHT edges
$ws.Range("A1:H1").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous
$ws.Range("H2:H12").Borders.Item($xlEdgeRight).LineStyle = $xlContinuous

# TOP edges
$ws.Range("A1:H1").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("A2:H2").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("H6").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous
$ws.Range("A12:H12").Borders.Item($xlEdgeTop).LineStyle = $xlContinuous

# BOTTOM edges
$ws.Range("A2").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("A5:G5").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("A6").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("A11:G11").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous
$ws.Range("A12:H12").Borders.Item($xlEdgeBottom).LineStyle = $xlContinuous

# 7) View state: scroll to column D and select H16 (best-effort, matches author's last saved selection)
$ws.Range("H16").Select()
